$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows scraped from MV (market data) for Aug/Sep 2021, continuing
# the "En $ 2021 - Diaria" series right after row 147 (03-08-2021).
$newRows = @(
    @{Row=148; A="04-08-2021"; C=2.35; D=3.67; E=4.46},
    @{Row=149; A="05-08-2021"; C=2.35; D=3.72; E=4.57},
    @{Row=150; A="06-08-2021"; C=2.35; D=3.78; E=4.64},
    @{Row=151; A="09-08-2021"; C=$null; D=3.85; E=4.67},
    @{Row=152; A="10-08-2021"; C=2.32; D=3.9;  E=4.67},
    @{Row=153; A="11-08-2021"; C=2.55; D=3.99; E=4.76},
    @{Row=154; A="12-08-2021"; C=2.73; D=4;    E=4.75},
    @{Row=155; A="13-08-2021"; C=$null; D=4;   E=4.76},
    @{Row=156; A="16-08-2021"; C=2.77; D=3.97; E=4.72},
    @{Row=157; A="17-08-2021"; C=$null; D=3.91; E=4.73},
    @{Row=158; A="18-08-2021"; C=2.8;  D=3.95; E=4.77},
    @{Row=159; A="19-08-2021"; C=2.86; D=3.98; E=4.85},
    @{Row=160; A="20-08-2021"; C=2.9;  D=4.01; E=4.87},
    @{Row=161; A="23-08-2021"; C=2.92; D=4.02; E=4.89},
    @{Row=162; A="24-08-2021"; C=$null; D=4.02; E=4.86},
    @{Row=163; A="25-08-2021"; C=$null; D=3.99; E=4.79},
    @{Row=164; A="26-08-2021"; C=$null; D=3.94; E=4.76},
    @{Row=165; A="27-08-2021"; C=$null; D=3.94; E=4.74},
    @{Row=166; A="30-08-2021"; C=3.01; D=$null; E=4.72},
    @{Row=167; A="31-08-2021"; C=3.01; D=3.97; E=4.74},
    @{Row=168; A="01-09-2021"; C=3.67; D=4.4;  E=4.98},
    @{Row=169; A="02-09-2021"; C=3.82; D=4.41; E=5.09}
)

foreach ($r in $newRows) {
    # Force the date label into column A as literal text (matching the rest
    # of the "Serie" column) instead of letting it auto-parse into a date
    # serial number. Pre-format as Text, assign, then drop back to the
    # workbook's "Normal" cell style so no stray number format lingers.
    $cell = $ws.Cells.Item($r.Row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $r.A
    $cell.Style = "Normal"

    if ($null -ne $r.C) { $ws.Cells.Item($r.Row, 3).Value = $r.C }
    if ($null -ne $r.D) { $ws.Cells.Item($r.Row, 4).Value = $r.D }
    if ($null -ne $r.E) { $ws.Cells.Item($r.Row, 5).Value = $r.E }
}
